# Update the "想去人数" (interested-people count) figures in column F
# for the "展览" and "全部类型" sheets, reflecting newly generated data
# (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 4238
$wsExhibit.Range("F3").Value = 2413
$wsExhibit.Range("F4").Value = 479
$wsExhibit.Range("F7").Value = 47
$wsExhibit.Range("F10").Value = 123
$wsExhibit.Range("F11").Value = 149
$wsExhibit.Range("F12").Value = 1580
$wsExhibit.Range("F14").Value = 3245
$wsExhibit.Range("F15").Value = 220

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 4238
$wsAll.Range("F3").Value = 2413
$wsAll.Range("F4").Value = 479
$wsAll.Range("F8").Value = 47
$wsAll.Range("F12").Value = 123
$wsAll.Range("F13").Value = 149
$wsAll.Range("F16").Value = 1580
$wsAll.Range("F18").Value = 3245
$wsAll.Range("F19").Value = 220
